# Commit before changing all of the names
#
# This script reproduces, via Excel COM interop, the changes observed in the
# canonical-OOXML diff:
#  - ~NMOS_GSD ("sheet4"): column G header "JLCPCB Type" -> "ki_keywords",
#    and its two data rows "Basic" -> "transistor NMOS N-MOS N-MOSFET basic".
#  - ~PMOS_GSD ("sheet10"): a new "ki_keywords" column is inserted between
#    the old "JLCPCB Type" (now "Description") and "Description" (now
#    "ki_keywords") columns; the "JLCPCB Type"/"Basic" column is repurposed
#    to hold the Description text, and a brand-new column holds the new
#    keyword strings.
#  - Column width / selection cosmetic updates on both sheets.
#  - Workbook window position/size + firstSheet bookview cosmetics.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ~NMOS_GSD: rename the "JLCPCB Type" column to "ki_keywords" and replace
# the "Basic" values with real keyword strings.
# ---------------------------------------------------------------------
$wsNmosGsd = $wb.Worksheets.Item("~NMOS_GSD")
$wsNmosGsd.Range("G1").Value = "ki_keywords"
$wsNmosGsd.Range("G2").Value = "transistor NMOS N-MOS N-MOSFET basic"
$wsNmosGsd.Range("G3").Value = "transistor NMOS N-MOS N-MOSFET basic"

# ---------------------------------------------------------------------
# ~PMOS_GSD: the old column D ("JLCPCB Type" / "Basic") becomes the
# Description column (header + values move in from the old column E), and
# the old column E becomes the new ki_keywords column.
# ---------------------------------------------------------------------
$wsPmosGsd = $wb.Worksheets.Item("~PMOS_GSD")
$wsPmosGsd.Range("D1").Value = "Description"
$wsPmosGsd.Range("D2").Value = "P-Channel 30V 4A 1.4W 85mΩ@2.5V SOT-23 MOSFET"
$wsPmosGsd.Range("D3").Value = "P-Channel 20V 2.3A 1.6W 142mΩ@2.5V SOT-23 MOSFET"

$wsPmosGsd.Range("E1").Value = "ki_keywords"
$wsPmosGsd.Range("E2").Value = "transistor PMOS P-MOS P-MOSFET basic"
$wsPmosGsd.Range("E3").Value = "transistor PMOS P-MOS P-MOSFET basic"

# ---------------------------------------------------------------------
# Column width cosmetics.
# ---------------------------------------------------------------------
$wsNmosGsd.Columns.Item(7).ColumnWidth = 8.43
$wsPmosGsd.Columns.Item(4).ColumnWidth = 50
$wsPmosGsd.Columns.Item(5).ColumnWidth = 50

# ---------------------------------------------------------------------
# Selections: set the inactive sheet's selection first, then the active
# sheet's selection last so it stays the active tab.
# ---------------------------------------------------------------------
$wsNmosGsd.Range("J10").Select()
$wsPmosGsd.Range("E2:E3").Select()

# ---------------------------------------------------------------------
# Workbook window geometry / first visible tab.
# ---------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = 3315
$win.Top = 1005
$win.Width = 24840
$win.Height = 12450

$wb.Worksheets.Item("~NMOS_SGD").Activate()
$wsPmosGsd.Activate()
